$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (D3) now carries the "API testing" question. The old "progarms"
# typo question that used to live in D3 is retired from the workbook.
$ws.Range("D3").Value = "Create API testing interview questions covering REST principles, HTTP methods, status codes, authentication, request/response validation, and Rest Assured.  Note: Always provide the question in serial number format"
$ws.Rows("3").RowHeight = 55.2

# Rows 4-6 keep their original questions (Selenium / TestNG / Maven) - no
# text change needed there; only the shared-string bookkeeping shifts.

# Rows 7-11 previously held the API testing / Include.../ Mix.../ Target...
# / Do not include... questions - that data is removed from the sheet.
$ws.Range("A7:D11").ClearContents()
$ws.Range("A7:D11").Style = "Normal"

# Selection moves to C4 to match the new cursor position.
$ws.Range("C4").Select()
